$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row text updates (lowercased / renamed) ---
$ws.Range("A1").Value = "industry"
$ws.Range("B1").Value = "unit"
$ws.Range("C1").Value = "process"
$ws.Range("D1").Value = "carbon (kg CO2 eq)"
$ws.Range("E1").Value = "ced (MJ)"
$ws.Range("F1").Value = "climate change (kg CO2 eq)"
$ws.Range("G1").Value = "region"

# --- Header comments describing each columns data type ---
$ws.Range("A1").AddComment("Data type: Categorical (text)") | Out-Null
$ws.Range("B1").AddComment("Data type: Various (e.g. kg, kWh)") | Out-Null
$ws.Range("C1").AddComment("Data type: Categorical (text)") | Out-Null
$ws.Range("D1").AddComment("Data type: Carbon footprint") | Out-Null
$ws.Range("E1").AddComment("Data type: Cumulative energy demand") | Out-Null
$ws.Range("F1").AddComment("Data type: Climate change impact") | Out-Null
$ws.Range("G1").AddComment("Data type: Categorical (text)") | Out-Null

# --- Data rows: D now holds old carbon(E) values, E now holds old CED(F) values, ---
# --- and F holds newly supplied climate-change-impact figures. ---
$ws.Range("D2").Value = 5.13
$ws.Range("E2").Value = 100.338
$ws.Range("F2").Value = 0.00014303801
$ws.Range("D3").Value = 0.3329423266666667
$ws.Range("E3").Value = 5.7885818
$ws.Range("F3").Value = 0.0000092833156
$ws.Range("D4").Value = 0.7626977333333333
$ws.Range("E4").Value = 13.23279
$ws.Range("F4").Value = 0.000021266037
$ws.Range("D5").Value = 0.8813826000000001
$ws.Range("E5").Value = 13.630893
$ws.Range("F5").Value = 0.000024575285
$ws.Range("D6").Value = 1.313870133333333
$ws.Range("E6").Value = 21.129109
$ws.Range("F6").Value = 0.000036634185
$ws.Range("D7").Value = 1.000428
$ws.Range("E7").Value = 15.678095
$ws.Range("F7").Value = 0.000027894588
$ws.Range("D8").Value = 3.295048733333334
$ws.Range("E8").Value = 107.43137
$ws.Range("F8").Value = 0.000091874702
$ws.Range("D9").Value = 7.504679333333333
$ws.Range("E9").Value = 115.2773
$ws.Range("F9").Value = 0.00020925037
$ws.Range("D10").Value = 3.797658933333334
$ws.Range("E10").Value = 41.685644
$ws.Range("F10").Value = 0.00010588881
$ws.Range("D11").Value = 3.797658933333334
$ws.Range("E11").Value = 41.685645
$ws.Range("F11").Value = 0.00010588881
$ws.Range("D12").Value = 3.9534742
$ws.Range("E12").Value = 49.063453
$ws.Range("F12").Value = 0.00011023335
$ws.Range("D13").Value = 2.742414666666667
$ws.Range("E13").Value = 43.724239
$ws.Range("F13").Value = 0.000076465798
$ws.Range("D14").Value = 1.332928466666667
$ws.Range("E14").Value = 22.840126
$ws.Range("F14").Value = 0.000037165582
$ws.Range("D15").Value = 5.46
$ws.Range("E15").Value = 124.2204
$ws.Range("F15").Value = 0.00015223929
$ws.Range("D16").Value = 2.051337733333333
$ws.Range("E16").Value = 72.04149099999999
$ws.Range("F16").Value = 0.00005719674
$ws.Range("D17").Value = 3.543347266666667
$ws.Range("E17").Value = 73.408489
$ws.Range("F17").Value = 0.00009879792299999999
$ws.Range("D18").Value = 3.324106266666667
$ws.Range("E18").Value = 74.167906
$ws.Range("F18").Value = 0.000092684902
$ws.Range("D19").Value = 1.5
$ws.Range("E19").Value = 21.783
$ws.Range("F19").Value = 0.00004182398
$ws.Range("D20").Value = 78.08035333333333
$ws.Range("E20").Value = 61.672334
$ws.Range("F20").Value = 0.0021770875
$ws.Range("D21").Value = 1.63
$ws.Range("E21").Value = 82.255
$ws.Range("F21").Value = 0.000045448725
$ws.Range("D22").Value = 1.377814533333333
$ws.Range("E22").Value = 10.828352
$ws.Range("F22").Value = 0.000038417125
$ws.Range("D23").Value = 3.115061133333334
$ws.Range("E23").Value = 42.761759
$ws.Range("F23").Value = 0.00008685616999999999
$ws.Range("D24").Value = 3.115061133333334
$ws.Range("E24").Value = 42.761759
$ws.Range("F24").Value = 0.00008685616999999999
$ws.Range("D25").Value = 1.130179133333333
$ws.Range("E25").Value = 16.810616
$ws.Range("F25").Value = 0.000031512393
$ws.Range("D26").Value = 0.2490375666666667
$ws.Range("E26").Value = 4.4964375
$ws.Range("F26").Value = 0.0000069438282
$ws.Range("D27").Value = 0.10254488
$ws.Range("E27").Value = 1.8514743
$ws.Range("F27").Value = 0.0000028592234
$ws.Range("D28").Value = 0.06836325333333335
$ws.Range("E28").Value = 1.2343162
$ws.Range("F28").Value = 0.0000019061489
$ws.Range("D29").Value = 1.7232184
$ws.Range("E29").Value = 23.923034
$ws.Range("F29").Value = 0.000048047902
$ws.Range("D30").Value = 10.95765333333333
$ws.Range("E30").Value = 197.84325
$ws.Range("F30").Value = 0.00030552844
$ws.Range("D31").Value = 7.046298000000001
$ws.Range("E31").Value = 127.22273
$ws.Range("F31").Value = 0.00019646949
$ws.Range("D32").Value = 4.931920466666667
$ws.Range("E32").Value = 89.047096
$ws.Range("F32").Value = 0.00013751503
$ws.Range("D33").Value = 3.286319266666667
$ws.Range("E33").Value = 59.335342
$ws.Range("F33").Value = 0.000091631302
$ws.Range("D34").Value = 2.465960266666666
$ws.Range("E34").Value = 44.523548
$ws.Range("F34").Value = 0.00006875751499999999
$ws.Range("D35").Value = 1.6456012
$ws.Range("E35").Value = 29.711754
$ws.Range("F35").Value = 0.000045883728
$ws.Range("D36").Value = 1.235421666666667
$ws.Range("E36").Value = 22.305857
$ws.Range("F36").Value = 0.000034446834
$ws.Range("D37").Value = 0.9863840666666667
$ws.Range("E37").Value = 17.809419
$ws.Range("F37").Value = 0.000027503006
$ws.Range("D38").Value = 0.7760302666666667
$ws.Range("E38").Value = 11.87686
$ws.Range("F38").Value = 0.000021637784
$ws.Range("D39").Value = 0.22582928
$ws.Range("E39").Value = 3.9291096
$ws.Range("F39").Value = 0.0000062967195
$ws.Range("D40").Value = 0.4883089600000001
$ws.Range("E40").Value = 8.8165441
$ws.Range("F40").Value = 0.000013615349
$ws.Range("D41").Value = 48.18144466666666
$ws.Range("E41").Value = 869.92841
$ws.Range("F41").Value = 0.0013434265
$ws.Range("D42").Value = 24.08828066666667
$ws.Range("E42").Value = 434.92012
$ws.Range("F42").Value = 0.00067164519
$ws.Range("D43").Value = 16.06048133333334
$ws.Range("E43").Value = 289.97614
$ws.Range("F43").Value = 0.00044780884
$ws.Range("D44").Value = 10.32285133333333
$ws.Range("E44").Value = 186.38174
$ws.Range("F44").Value = 0.00028782849
$ws.Range("D45").Value = 7.226972666666667
$ws.Range("E45").Value = 130.48485
$ws.Range("F45").Value = 0.00020150717
$ws.Range("D46").Value = 4.8196094
$ws.Range("E46").Value = 87.01929
$ws.Range("F46").Value = 0.0001343835
$ws.Range("D47").Value = 3.613486266666667
$ws.Range("E47").Value = 65.24242599999999
$ws.Range("F47").Value = 0.00010075359
$ws.Range("D48").Value = 2.4073632
$ws.Range("E48").Value = 43.465562
$ws.Range("F48").Value = 0.000067123673
$ws.Range("D49").Value = 1.806743133333333
$ws.Range("E49").Value = 32.621213
$ws.Range("F49").Value = 0.000050376793
$ws.Range("D50").Value = 1.445394533333333
$ws.Range("E50").Value = 26.096971
$ws.Range("F50").Value = 0.000040301434
$ws.Range("D51").Value = 3.8514248
$ws.Range("E51").Value = 120.22883
$ws.Range("F51").Value = 0.00010738794
$ws.Range("D52").Value = 0.4297642400000001
$ws.Range("E52").Value = 6.0446765
$ws.Range("F52").Value = 0.000011982967
$ws.Range("D53").Value = 0.31720566
$ws.Range("E53").Value = 5.3673553
$ws.Range("F53").Value = 0.0000088445355
$ws.Range("D54").Value = 3.315677
$ws.Range("E54").Value = 92.73489499999999
$ws.Range("F54").Value = 0.000092449873
